$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.0002446666666666666
$ws.Range("H2").Value = 0.000734
$ws.Range("I2").Value = 0.001827111446330468
$ws.Range("J2").Value = 0.001827111446330468
$ws.Range("M2").Value = 721.1356606666667
$ws.Range("N2").Value = 2163.406982
$ws.Range("O2").Value = 0.8508208584329936
$ws.Range("P2").Value = 0.8508208584329936
$ws.Range("Q2").Value = 0.1764378583097778
$ws.Range("R2").Value = 1.587940724788
$ws.Range("S2").Value = 0.001554544529219638
$ws.Range("T2").Value = 0.001554544529219638

# Row 3
$ws.Range("G3").Value = 0.0002446666666666666
$ws.Range("H3").Value = 0.000734
$ws.Range("I3").Value = 0.001827111446330468
$ws.Range("J3").Value = 0.001827111446330468
$ws.Range("O3").Value = 0.002793596814304166
$ws.Range("P3").Value = 0.002793596814304166
$ws.Range("Q3").Value = 0.0005793184711111111
$ws.Range("R3").Value = 0.00521386624
$ws.Range("S3").Value = 0.000005104212715847474
$ws.Range("T3").Value = 0.000005104212715847474

# Row 4
$ws.Range("G4").Value = 0.0002446666666666666
$ws.Range("H4").Value = 0.000734
$ws.Range("I4").Value = 0.001827111446330468
$ws.Range("J4").Value = 0.001827111446330468
$ws.Range("O4").Value = 0.1463855447527022
$ws.Range("P4").Value = 0.1463855447527022
$ws.Range("Q4").Value = 0.03035651012511111
$ws.Range("R4").Value = 0.273208591126
$ws.Range("S4").Value = 0.0002674627043949831
$ws.Range("T4").Value = 0.0002674627043949832

# Row 5
$ws.Range("I5").Value = 0.9981728885536695
$ws.Range("J5").Value = 0.9981728885536695
$ws.Range("M5").Value = 721.1356606666667
$ws.Range("N5").Value = 2163.406982
$ws.Range("O5").Value = 0.8508208584329936
$ws.Range("P5").Value = 0.8508208584329936
$ws.Range("Q5").Value = 96.39011732590286
$ws.Range("R5").Value = 867.5110559331258
$ws.Range("S5").Value = 0.8492663139037739
$ws.Range("T5").Value = 0.8492663139037739

# Row 6
$ws.Range("I6").Value = 0.9981728885536695
$ws.Range("J6").Value = 0.9981728885536695
$ws.Range("O6").Value = 0.002793596814304166
$ws.Range("P6").Value = 0.002793596814304166
$ws.Range("S6").Value = 0.002788492601588318
$ws.Range("T6").Value = 0.002788492601588318

# Row 7
$ws.Range("I7").Value = 0.9981728885536695
$ws.Range("J7").Value = 0.9981728885536695
$ws.Range("O7").Value = 0.1463855447527022
$ws.Range("P7").Value = 0.1463855447527022
$ws.Range("S7").Value = 0.1461180820483072
$ws.Range("T7").Value = 0.1461180820483072
